$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (sheet 1): refresh "want to go" counts (col F) and add the new
# "排球少年" event as a row inserted just before the trailing "银魂" row.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 524
$ws1.Range("F3").Value = 777
$ws1.Range("F4").Value = 1562
$ws1.Range("F5").Value = 240
$ws1.Range("F7").Value = 179
$ws1.Range("F9").Value = 6366
$ws1.Range("F11").Value = 413
$ws1.Range("F13").Value = 5525
$ws1.Range("F15").Value = 183
$ws1.Range("F16").Value = 1210
$ws1.Range("F17").Value = 4
$ws1.Range("F18").Value = 67
$ws1.Range("F19").Value = 369
$ws1.Range("F20").Value = 76
$ws1.Range("F22").Value = 315
$ws1.Range("F24").Value = 1
$ws1.Range("F25").Value = 3954

# Insert a fresh row 26, pushing the old row 26 ("银魂主题派对only2.0") down
# to row 27, then copy the formatting of the row above onto the new index
# cell so it keeps the shared "s=1" look.
$ws1.Rows.Item(26).Insert()
$ws1.Range("A25").Copy()
$ws1.Range("A26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A26").Value = 25
# Force the start-date text to stay a plain string (it looks like a date, so
# a bare .Value assignment would get auto-coerced into a real date value) -
# mark the cell as text first, then restore the plain "General" look from a
# neighbouring text cell once the literal string has been stored.
$ws1.Range("B26").NumberFormat = "@"
$ws1.Range("B26").Value = "2024-08-10"
$ws1.Range("B25").Copy()
$ws1.Range("B26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("C26").Value = "合肥·排球少年only之夏日招新季"
$ws1.Range("D26").Value = "广德路与长江东路交口往北200米文一时埠里文旅街区 巅峰篮球公园"
$ws1.Range("E26").Value = "2024.08.10 10:00-08.10 17:00"
$ws1.Range("F26").Value = 7
$ws1.Range("G26").Value = 70
$ws1.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=88281"
$ws1.Range("I26").Value = "//i0.hdslb.com/bfs/openplatform/202406/qjd7yzXE1719556597555.jpeg"

$ws1.Range("A27").Value = 26
$ws1.Range("F27").Value = 169

# ---------------------------------------------------------------------------
# Sheet "演出" (sheet 2): refresh the single "want to go" count.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 100

# ---------------------------------------------------------------------------
# Sheet "全部类型" (sheet 4): same refresh as above, offset by one row
# because row 2 holds the "演出" entry ahead of the "展览" ones.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 100
$ws4.Range("F3").Value = 524
$ws4.Range("F4").Value = 777
$ws4.Range("F5").Value = 1562
$ws4.Range("F6").Value = 240
$ws4.Range("F8").Value = 179
$ws4.Range("F10").Value = 6366
$ws4.Range("F12").Value = 413
$ws4.Range("F14").Value = 5525
$ws4.Range("F16").Value = 183
$ws4.Range("F17").Value = 1210
$ws4.Range("F18").Value = 4
$ws4.Range("F19").Value = 67
$ws4.Range("F20").Value = 369
$ws4.Range("F21").Value = 76
$ws4.Range("F23").Value = 315
$ws4.Range("F25").Value = 1
$ws4.Range("F26").Value = 3954

$ws4.Rows.Item(28).Insert()
$ws4.Range("A27").Copy()
$ws4.Range("A28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws4.Range("A28").Value = 27
$ws4.Range("B28").NumberFormat = "@"
$ws4.Range("B28").Value = "2024-08-10"
$ws4.Range("B27").Copy()
$ws4.Range("B28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws4.Range("C28").Value = "合肥·排球少年only之夏日招新季"
$ws4.Range("D28").Value = "广德路与长江东路交口往北200米文一时埠里文旅街区 巅峰篮球公园"
$ws4.Range("E28").Value = "2024.08.10 10:00-08.10 17:00"
$ws4.Range("F28").Value = 7
$ws4.Range("G28").Value = 70
$ws4.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=88281"
$ws4.Range("I28").Value = "//i0.hdslb.com/bfs/openplatform/202406/qjd7yzXE1719556597555.jpeg"

$ws4.Range("A29").Value = 28
$ws4.Range("F29").Value = 169
